$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.129.39'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.43%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.783.42'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.63%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.36%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '336.16'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.20%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3830'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.14%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3429'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.55%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.25'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.74%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.195'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.04%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07501'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.33%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.005'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.22%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.79'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.88%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.456'

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.787.51'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.41%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.086'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.70%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001093'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.08%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06690'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.43%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '83.78'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.13%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.003'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.25%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.634'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.59%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.35'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.34%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.153.15'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.36%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.32'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -6.28%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.378'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.91%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.535'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -5.37%  '

# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.27'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.04%  '

# Row 28
$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.469'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.53%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '154.29'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.28%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.990.57'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.38%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '134.64'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.51%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.018'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.63%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.060'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -4.99%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08711'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.29%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.24'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -5.04%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.653'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -4.08%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6895'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.03%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.422'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.92%  '

# Row 39
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06333'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.31%  '

# Row 40
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2194'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.20%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.778'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.72%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.02338'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.95%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.238'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.45%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.39'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.65%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6481'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.29%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.003'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.34%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.852'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -4.65%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.142'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.91%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '130.04'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.37%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07126'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.20%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.12'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.02%  '
